$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 9.8
$ws.Columns.Item(3).ColumnWidth = 11.8

$ws.Range("A1").Value = 161.11021704595609
$ws.Range("B1").Value = 6.6404327169427013
$ws.Range("C1").Value = 0.55714285714285716
